# updated signs on MFA data
# Flip the sign of the flux values (columns G:N) for the nutrient-uptake
# exchange reactions on the "Rxns" sheet, rows 105,107,109,110,113-120.
# Also restore the last-saved selection state recorded in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rxns")
$ws.Activate()

$rows = @(105, 107, 109, 110, 113, 114, 115, 116, 117, 118, 119, 120)
$cols = @("G", "H", "I", "J", "K", "L", "M", "N")

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $cell.Value2 = -1 * $cell.Value2
    }
}

# Restore the recorded selection (activeCell G113) from the saved workbook view.
$ws.Range("G113").Select() | Out-Null
